$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column P/Q data: 9-mer peptide + gj score, one row per existing
# data row (rows 4-164), matching the sliding-window layout already
# used by the other peptide/score column pairs on this sheet.
$pqData = @(
    @(4, "MEEFEITIH", 0.06572004464283529),
    @(5, "MMEEFEITI", 0.06527797505648819),
    @(6, "EEFEITIHR", 0.06380417009806016),
    @(7, "EMMEEFEIT", 0.06364996497479281),
    @(8, "LLGYIEEIK", 0.0622395860008393),
    @(9, "EFEITIHRP", 0.06178846466875153),
    @(10, "LGYIEEIKF", 0.06107992215464007),
    @(11, "GYIEEIKFA", 0.05960548647190239),
    @(12, "IAGAANWTN", 0.05923832575224897),
    @(13, "AGAANWTNG", 0.05922711262346664),
    @(14, "GKIHILAFK", 0.05835021713324512),
    @(15, "GLLGYIEEI", 0.05821989058040858),
    @(16, "FEITIHRPK", 0.05745634206169023),
    @(17, "DGKIHILAF", 0.05717046465673384),
    @(18, "YIEEIKFAY", 0.05703672622349282),
    @(19, "GIAGAANWT", 0.05688675962120156),
    @(20, "ESGIAGAAN", 0.0568722155900105),
    @(21, "AAARGTITL", 0.05673188251400917),
    @(22, "AARGTITLT", 0.05579545683740722),
    @(23, "AAAARGTIT", 0.05567950941874505),
    @(24, "ARGTITLTK", 0.05545172428175964),
    @(25, "DLTELLFTY", 0.05502042417226133),
    @(26, "AESGIAGAA", 0.05484182214837124),
    @(27, "REMMEEFEI", 0.05442942614185246),
    @(28, "KIHILAFKN", 0.054398745973054),
    @(29, "EGLLGYIEE", 0.05406531269718829),
    @(30, "IEEIKFAYS", 0.05391825530687247),
    @(31, "GDLTELLFT", 0.05375791813564762),
    @(32, "YSLEHAESG", 0.05362787166492214),
    @(33, "EHAESGIAG", 0.05341597290906665),
    @(34, "TPAAAARGT", 0.05333119690076816),
    @(35, "PAAAARGTI", 0.05332114644339227),
    @(36, "SGIAGAANW", 0.0532573723988755),
    @(37, "GYSLEHAES", 0.05317635366526818),
    @(38, "FEKVLITHM", 0.05304344777516209),
    @(39, "VVQEGHDGK", 0.05304340963382105),
    @(40, "LTPAAAARG", 0.05276560147003617),
    @(41, "HDGKIHILA", 0.0526819631475348),
    @(42, "LTELLFTYK", 0.05267431343811571),
    @(43, "LEHAESGIA", 0.0524016233403759),
    @(44, "GHDGKIHIL", 0.05240073753363172),
    @(45, "DSNGIKEGL", 0.05225886368830376),
    @(46, "EITIHRPKT", 0.0521463670926493),
    @(47, "PKTDTTGGD", 0.05213246486801888),
    @(48, "KTDTTGGDL", 0.05212624755155061),
    @(49, "KTQGEIKGS", 0.05212256497689984),
    @(50, "TTGGDLTEL", 0.0520317326370146),
    @(51, "SNGIKEGLL", 0.05201514018893391),
    @(52, "HAESGIAGA", 0.05181961920190777),
    @(53, "TGGDLTELL", 0.05179597659350441),
    @(54, "RGTITLTKE", 0.05175292829540754),
    @(55, "EKVLITHMD", 0.05174443887591526),
    @(56, "YKFEKVLIT", 0.0513844918242553),
    @(57, "GLTPAAAAR", 0.05132097010488448),
    @(58, "RLQEGLTPA", 0.05128662152718959),
    @(59, "SLEHAESGI", 0.05118052802081148),
    @(60, "KVLITHMDQ", 0.05115509530012802),
    @(61, "TDTTGGDLT", 0.0511509484222384),
    @(62, "GGDLTELLF", 0.05101773968427599),
    @(63, "ELLFTYKFE", 0.05089841314185004),
    @(64, "KGKTQGEIK", 0.05087163342982749),
    @(65, "LLFTYKFEK", 0.05085577143521737),
    @(66, "EGLTPAAAA", 0.05081461334159701),
    @(67, "GSVVQEGHD", 0.05081340777503747),
    @(68, "DTTGGDLTE", 0.05064095776793608),
    @(69, "LQEGLTPAA", 0.05054296987924885),
    @(70, "KEGLLGYIE", 0.05054210270650397),
    @(71, "GKTQGEIKG", 0.05036243234915019),
    @(72, "RPKTDTTGG", 0.05023487554452297),
    @(73, "MLAGIYLKV", 0.05023335775740555),
    @(74, "SVVQEGHDG", 0.05023068113096193),
    @(75, "GIKEGLLGY", 0.05012476252685333),
    @(76, "TQGEIKGSV", 0.05003250019529507),
    @(77, "TYKFEKVLI", 0.0497364164940514),
    @(78, "QEGLTPAAA", 0.04964483531199283),
    @(79, "MPARLQEGL", 0.04962072190563029),
    @(80, "EGHDGKIHI", 0.04942560411861907),
    @(81, "IHILAFKND", 0.0494129134140244),
    @(82, "DDSNGIKEG", 0.04937057248782457),
    @(83, "IKEGLLGYI", 0.04936558150647609),
    @(84, "KFEKVLITH", 0.04930373840721156),
    @(85, "KGSVVQEGH", 0.04917547280310664),
    @(86, "DMPARLQEG", 0.04901001598359366),
    @(87, "TELLFTYKF", 0.04896789469731195),
    @(88, "VQEGHDGKI", 0.04895220676256098),
    @(89, "QEGHDGKIH", 0.048936235318655),
    @(90, "KREMMEEFE", 0.04891559670887716),
    @(91, "ARLQEGLTP", 0.04876689027392193),
    @(92, "GTITLTKEM", 0.04873492687604183),
    @(93, "LAGIYLKVK", 0.04866969260615035),
    @(94, "GKREMMEEF", 0.04833732372850617),
    @(95, "PARLQEGLT", 0.04829101766532597),
    @(96, "EIKFAYSGY", 0.0481120431354161),
    @(97, "SPQFLQALG", 0.04799156836793758),
    @(98, "HILAFKNDY", 0.04798684532876932),
    @(99, "NGIKEGLLG", 0.04773293111605311),
    @(100, "HKDDSNGIK", 0.04760569579498208),
    @(101, "ITIHRPKTD", 0.04756802751254531),
    @(102, "SGYSLEHAE", 0.04749208994787968),
    @(103, "ILAFKNDYD", 0.04726952769060278),
    @(104, "LQALGKREM", 0.04698079926815035),
    @(105, "QGEIKGSVV", 0.04692309443860947),
    @(106, "LFTYKFEKV", 0.04678094256787046),
    @(107, "EEIKFAYSG", 0.0466524042290422),
    @(108, "YSGYSLEHA", 0.04648246124556583),
    @(109, "PQFLQALGK", 0.04646474084183556),
    @(110, "IKFAYSGYS", 0.04646357574431557),
    @(111, "FTYKFEKVL", 0.04631780125602104),
    @(112, "YDMPARLQE", 0.04586157046510425),
    @(113, "SSPQFLQAL", 0.0457632991721247),
    @(114, "VLITHMDQY", 0.04560893748245581),
    @(115, "HRPKTDTTG", 0.04546453880074268),
    @(116, "TIHRPKTDT", 0.04540323165577202),
    @(117, "QALGKREMM", 0.04529440362994468),
    @(118, "KDDSNGIKE", 0.04528918345314785),
    @(119, "IKGSVVQEG", 0.04525827772788377),
    @(120, "LGKREMMEE", 0.04515044681294874),
    @(121, "NDYDMPARL", 0.0450773907189657),
    @(122, "DYDMPARLQ", 0.04485515090829116),
    @(123, "YSPTPHKDD", 0.0448389168486306),
    @(124, "IHRPKTDTT", 0.04473156556997555),
    @(125, "FLQALGKRE", 0.04463301604647513),
    @(126, "LTKEMDRSS", 0.04453556344412189),
    @(127, "AYSGYSLEH", 0.04445657749934873),
    @(128, "SPTPHKDDS", 0.04420300093786664),
    @(129, "RSSPQFLQA", 0.04418316447818018),
    @(130, "FKNDYDMPA", 0.04411445151791293),
    @(131, "TITLTKEMD", 0.04406887615833776),
    @(132, "LAFKNDYDM", 0.04398837561251283),
    @(133, "QYSPTPHKD", 0.04398080381280432),
    @(134, "KNDYDMPAR", 0.04396975525805932),
    @(135, "DQYSPTPHK", 0.04376608378560849),
    @(136, "VKGKTQGEI", 0.04361218671343957),
    @(137, "AFKNDYDMP", 0.04358079566780454),
    @(138, "TLTKEMDRS", 0.04343680118369748),
    @(139, "ITLTKEMDR", 0.04343648794603256),
    @(140, "PHKDDSNGI", 0.04315487687649518),
    @(141, "EIKGSVVQE", 0.0429368338646928),
    @(142, "GEIKGSVVQ", 0.04267558189027705),
    @(143, "QFLQALGKR", 0.04247584612159816),
    @(144, "MDQYSPTPH", 0.04212895734196687),
    @(145, "LITHMDQYS", 0.04202576688413191),
    @(146, "PTPHKDDSN", 0.04172096084857192),
    @(147, "DRSSPQFLQ", 0.04158958702340361),
    @(148, "KFAYSGYSL", 0.04155851265711533),
    @(149, "AGIYLKVKG", 0.041538472119402),
    @(150, "ALGKREMME", 0.04122201205667409),
    @(151, "KVKGKTQGE", 0.04065163881207413),
    @(152, "HMDQYSPTP", 0.04026710684926018),
    @(153, "KEMDRSSPQ", 0.04021749284390257),
    @(154, "TKEMDRSSP", 0.04010905109276046),
    @(155, "TPHKDDSNG", 0.03981527750617909),
    @(156, "FAYSGYSLE", 0.03980693344526509),
    @(157, "YLKVKGKTQ", 0.03925731855386381),
    @(158, "GIYLKVKGK", 0.03920404601284776),
    @(159, "THMDQYSPT", 0.03918376409960991),
    @(160, "MDRSSPQFL", 0.03881034176952239),
    @(161, "EMDRSSPQF", 0.03871670275234826),
    @(162, "ITHMDQYSP", 0.03795686873921085),
    @(163, "LKVKGKTQG", 0.03781645974983729),
    @(164, "IYLKVKGKT", 0.03623898512546719)
)

foreach ($entry in $pqData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 16).Value = $entry[1]   # column P
    $ws.Cells.Item($r, 17).Value = $entry[2]   # column Q
}

